# Fruta / hortaliza, semanal
#
# Insert two new weekly records for "Ají" right before the existing row 688
# (pushing the previous rows 688-705 down to 690-707), then populate the
# two freshly inserted rows with the new data points.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 688 - everything currently at/after row
# 688 (through 705) shifts down to 690-707, and the sheet's used range
# grows from A1:R705 to A1:R707.
$ws.Rows.Item(688).Resize(2).Insert()

# New row 688: Ají, Americana (o), Primera
$ws.Range("A688").Value = 6
$ws.Range("B688").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C688").Value = "Metropolitana"
$ws.Range("D688").Value = 44595
$ws.Range("E688").Value = 13
$ws.Range("F688").Value = 100112021
$ws.Range("G688").Value = "Ají"
$ws.Range("H688").Value = "Americana (o)"
$ws.Range("I688").Value = "Primera"
$ws.Range("J688").Value = 140
$ws.Range("K688").Value = 20000
$ws.Range("L688").Value = 23000
$ws.Range("M688").Value = 21714
$ws.Range("N688").Value = "$/caja 25 kilos"
$ws.Range("O688").Value = "Provincia de Limarí"
$ws.Range("P688").Value = 869
$ws.Range("Q688").Value = 25
$ws.Range("R688").Value = "Hortaliza"

# New row 689: Ají, Chilena(o), Primera
$ws.Range("A689").Value = 6
$ws.Range("B689").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C689").Value = "Metropolitana"
$ws.Range("D689").Value = 44595
$ws.Range("E689").Value = 13
$ws.Range("F689").Value = 100112021
$ws.Range("G689").Value = "Ají"
$ws.Range("H689").Value = "Chilena(o)"
$ws.Range("I689").Value = "Primera"
$ws.Range("J689").Value = 140
$ws.Range("K689").Value = 27000
$ws.Range("L689").Value = 30000
$ws.Range("M689").Value = 28714
$ws.Range("N689").Value = "$/caja 25 kilos"
$ws.Range("O689").Value = "Provincia de Huasco"
$ws.Range("P689").Value = 1149
$ws.Range("Q689").Value = 25
$ws.Range("R689").Value = "Hortaliza"
